$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Netherlands: full repeater list (same shape as Germany/Portugal/etc.),
# based on "Portugal" sheet which already has matching column widths.
# ---------------------------------------------------------------------------
$srcNL = $wb.Worksheets.Item("Portugal")
$srcNL.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2200"
$nl.Range("B2").Value = "Netherlands Market"
[void]$nl.Range("B4").Select()

# ---------------------------------------------------------------------------
# Austria: repeater list without P32AR/P32DR, based on "Slovakia" sheet
# which already has that reduced shape + matching column widths.
# ---------------------------------------------------------------------------
$srcAT = $wb.Worksheets.Item("Slovakia")
$srcAT.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2307"
$at.Range("B2").Value = "Austria Market"
[void]$at.Range("B4").Select()

# ---------------------------------------------------------------------------
# Denmark: same reduced shape as Austria, plus an extra "MZXSDR240" row
# inserted right after "MZXDR240".
# ---------------------------------------------------------------------------
$srcDK = $wb.Worksheets.Item("Slovakia")
$srcDK.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"

$dk.Rows.Item(16).Insert()
$dk.Range("A15").Copy()
$dk.Range("A16").PasteSpecial(-4122)
$dk.Range("A16").Value = "MZXSDR240"

$dk.Range("B4").Value = "NGC-2913/T2796"
$dk.Range("B2").Value = "Denmark Market"

[void]$dk.Activate()
[void]$dk.Range("B15").Select()

# Cosmetic: try to match the tab-scroll position recorded in the target
# workbook (firstSheet). Not all engines expose/back this property, so
# failures here are swallowed.
try {
    $win = $excel.ActiveWindow
    $win.FirstSheet = 3
} catch {}
